$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (Days, Total cases, Daily cases, Deaths) for sheet rows 445-470,
# continuing the existing COVID tracking table.
$newRows = @(
    "444,663952,3978,10655",
    "445,665585,1633,10667",
    "446,667380,1795,10685",
    "447,669067,1687,10695",
    "448,670613,1546,10702",
    "449,670613,0,10704",
    "450,670613,0,10706",
    "451,674296,3683,10715",
    "452,675671,1375,10721",
    "453,677210,1539,10725",
    "454,677210,0,10728",
    "455,679510,2300,10730",
    "456,679510,0,10730",
    "457,679510,0,10731",
    "458,682160,2650,10744",
    "459,683400,1240,10751",
    "460,684954,1554,10755",
    "461,686152,1198,10760",
    "462,687353,1201,10770",
    "463,687353,0,10775",
    "464,687353,0,10776",
    "465,687353,0,10779",
    "466,690123,2770,10786",
    "467,691119,996,10789",
    "468,692111,992,10798",
    "469,693023,912,10801"
)

$startRow = 445
$r = $startRow
foreach ($line in $newRows) {
    $parts = $line.Split(",")
    $ws.Cells.Item($r, 1).Value = [int]$parts[0]
    $ws.Cells.Item($r, 2).Value = [double]$parts[1]
    $ws.Cells.Item($r, 3).Value = [double]$parts[2]
    $ws.Cells.Item($r, 4).Value = [double]$parts[3]
    $r = $r + 1
}
$endRow = $r - 1

# Copy the existing cell formatting (fills) down into the new rows so the
# new cells pick up the same style indices as the rest of the table.
$ws.Range("B444").Copy()
$ws.Range("B" + $startRow + ":B" + $endRow).PasteSpecial(-4122)

$ws.Range("C444").Copy()
$ws.Range("C" + $startRow + ":C" + $endRow).PasteSpecial(-4122)

$ws.Range("D2").Copy()
$ws.Range("D" + $startRow + ":D" + $endRow).PasteSpecial(-4122)

# Update the sheet selection to match the post-edit state (D2:D470 selected,
# active cell D2).
$ws.Range("D2:D470").Select()
